$d = $word.ActiveDocument

$bodyInner = @'
<w:p><w:r><w:t xml:space="preserve">O </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dongodb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> é um banco de dados de supermercado</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Data base </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dongodb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Conection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>localhost</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>User</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: admin</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Password</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>A10a20a30a40</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>$</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Para iniciar ele no </w:t></w:r><w:r><w:t xml:space="preserve">Visual Studio </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Code</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/>
'@

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Paragraphs.Item(1).Range
$target.Collapse(1)
$target.InsertXML($xml)
